$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" file. Mark it ready for handoff and
# bump its generate-date timestamp.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 12:41:48"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 ("b.md") picks up a new handoff file/date, its
# duplicate flag flips back to False, status becomes "Ready for handoff",
# and an error detail about a stale handback file is recorded. The Error
# Detail column also needs to be widened to fit the long message.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# A bare "False" is auto-coerced to a Boolean by the COM value-setter, but the
# source column is text ("True"/"False" strings); a leading apostrophe forces
# text, then re-copying a plain sibling's style drops the quote-prefix flag.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = $wsZhCn.Range("F2").Style
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-01 12:41:43"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/643dbef8e24dfdef65392a6975d144e216fb6aee/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/093163030e88e3d25082d03dd1c066ba22b035bd/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet: mirror of the zh-cn change above, for the de-de handoff file.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = $wsDeDe.Range("F2").Style
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-01 12:41:48"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/643dbef8e24dfdef65392a6975d144e216fb6aee/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/093163030e88e3d25082d03dd1c066ba22b035bd/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
